# Add "IDA for engines" dashboard-prep columns (ICAO Target / ATAG Target)
# and a header label for the (until now unlabeled) index column A.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -------------------------------------------------------
# A1 needs the same bold/bordered header style as B1:E1 already carry,
# so copy formatting from an existing header cell rather than re-building
# the style by hand (keeps the same cellXf / doesn't fork a new style).
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A1").Value = "Unnamed: 0"

$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "ICAO Target"

$ws.Range("G1").PasteSpecial(-4122)
$ws.Range("G1").Value = "ATAG Target"

$excel.CutCopyMode = 0

# --- Column A (rows 2:55) loses its bold/bordered style ---------------
# It now carries plain (unstyled) numeric formatting, matching the data
# columns B:E, so pull that formatting from B2 (an unstyled numeric cell).
$ws.Range("B2").Copy()
$ws.Range("A2:A55").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- New ICAO / ATAG target columns ------------------------------------
$icaoTargets = @{
  39 = 1.868127712959109
  40 = 1.830765158699927
  41 = 1.794149855525928
  42 = 1.75826685841541
  43 = 1.723101521247101
  44 = 1.688639490822159
  45 = 1.654866701005716
  46 = 1.621769366985602
  47 = 1.58933397964589
  48 = 1.557547300052972
  49 = 1.526396354051912
  50 = 1.495868426970874
  51 = 1.465951058431457
  52 = 1.436632037262827
  53 = 1.407899396517571
  54 = 1.379741408587219
  55 = 1.352146580415475
}

$atagTargets = @{
  45 = 1.641928483761755
  46 = 1.584460986830093
  47 = 1.52900485229104
  48 = 1.475489682460853
  49 = 1.423847543574724
  50 = 1.374012879549608
  51 = 1.325922428765372
  52 = 1.279515143758584
  53 = 1.234732113727033
  54 = 1.191516489746587
  55 = 1.149813412605456
}

for ($r = 2; $r -le 55; $r++) {
  if ($icaoTargets.ContainsKey($r)) {
    $ws.Cells.Item($r, 6).Value = $icaoTargets[$r]
  }
  if ($atagTargets.ContainsKey($r)) {
    $ws.Cells.Item($r, 7).Value = $atagTargets[$r]
  }
}
